$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held "1285870 - Marcos Villela Barcza" in B/C (no A label)
# is removed; everything below shifts up by one row.
$ws.Rows(13).Delete()

# After the shift, update the B/C (value) cells that now hold stale/incorrect
# text so they match the target content.
$ws.Range("B10:C10").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("B13:C13").Value = "Semestral"
$ws.Range("B15").Formula = "=""01/01/2021"""
$ws.Range("C15").Formula = "=""01/01/2021"""
$ws.Range("B15:C15").Copy()
$ws.Range("B15:C15").PasteSpecial(-4163)
$ws.Range("B18:C18").Value = "1285870 - Marcos Villela Barcza"
$ws.Range("B19:C19").Value = "Aulas expositivas, filmes e leituras de artigos técnicos"
$ws.Range("B20:C20").Value = "Provas e/ou trabalhos."
$ws.Range("B21:C21").Value = "Prova escrita para alunos que tenham média final maior ou igual a 3,0 (Três) e inferior a 5,0 (Cinco). A nota final será a média aritmética entre a média final e a prova escrita."
